$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some Price values in column D are plain numeric-looking strings (e.g. "302.74").
# Force those specific cells to Text format first so Excel keeps them stored as
# strings instead of silently converting them to numbers (matches the source data,
# which uses "." as a thousands-style separator rather than a decimal point).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "41.211.36"
$ws.Range("E2").Value = "  +2.45%  "
$ws.Range("D3").Value = "2.253.40"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "302.74"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").Value = "91.81"
$ws.Range("E6").Value = "  +3.39%  "
$ws.Range("D7").Value = "0.521"
$ws.Range("E7").Value = "  +1.68%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.483"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "32.21"
$ws.Range("E10").Value = "  +5.87%  "
$ws.Range("D11").Value = "52.66"
$ws.Range("E11").Value = "  +5.76%  "
$ws.Range("D12").Value = "0.0794"
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("E13").Value = "  +2.73%  "
$ws.Range("D14").Value = "6.59"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "2.602.37"
$ws.Range("E15").Value = "  +0.79%  "
$ws.Range("E16").Value = "  +1.80%  "
$ws.Range("D17").Value = "2.228.89"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "0.751"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("D19").Value = "41.124.51"
$ws.Range("E19").Value = "  +2.40%  "
$ws.Range("D20").Value = "11.95"
$ws.Range("E20").Value = "  +4.26%  "
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("E22").Value = "  +0.42%  "
$ws.Range("D23").Value = "66.74"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "239.76"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").Value = "2.56"
$ws.Range("E25").Value = "  +3.36%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").Value = "1.87"
$ws.Range("E27").Value = "  +2.80%  "
$ws.Range("D28").Value = "23.96"
$ws.Range("E28").Value = "  +5.15%  "
$ws.Range("D29").Value = "2.17"
$ws.Range("E29").Value = "  -1.74%  "
$ws.Range("D30").Value = "9.54"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "158.74"
$ws.Range("E31").Value = "  +2.31%  "
$ws.Range("D32").Value = "33.51"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("D34").Value = "5.12"
$ws.Range("E34").Value = "  +3.25%  "
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("D36").Value = "0.0731"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("E37").Value = "  +6.84%  "
$ws.Range("E38").Value = "  -0.50%  "
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").Value = "16.36"
$ws.Range("E40").Value = "  +3.37%  "
$ws.Range("D41").Value = "1.78"
$ws.Range("E41").Value = "  +5.21%  "
$ws.Range("D42").Value = "3.92"
$ws.Range("E42").Value = "  +1.99%  "
$ws.Range("D43").Value = "2.082.70"
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("D44").Value = "20.06"
$ws.Range("E44").Value = "  +10.62%  "
$ws.Range("D45").Value = "10.38"
$ws.Range("E45").Value = "  +6.01%  "
$ws.Range("D46").Value = "0.0276"
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("E47").Value = "  +8.56%  "
$ws.Range("D48").Value = "1.86"
$ws.Range("E48").Value = "  -12.81%  "
$ws.Range("D49").Value = "1.53"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").Value = "2.473.62"
$ws.Range("E50").Value = "  +1.09%  "
$ws.Range("D51").Value = "1.14"
$ws.Range("E51").Value = "  +4.17%  "
